$wb = $excel.ActiveWorkbook

# Update the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9122
$ws1.Range("F4").Value = 472

# Update the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9122
$ws4.Range("F4").Value = 472
